$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new backlog entry: category "Cleanup" with task "Remove code that's commented out"
$ws.Range("A16").Value = "Cleanup"
$ws.Range("B16").Value = "Remove code that's commented out"

# Move the active selection to A17, as in the saved workbook after the edit
$ws.Range("A17").Select()
